# Auto-generated Excel COM-interop script applying the cryptos.xlsx update
# Commit: Updated cryptos list on Tue Aug 29 02:52:05 UTC 2023 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.185.39'
$ws.Range("E2").Value = '  -0.05%  '

$ws.Range("D3").Value = '1.657.54'
$ws.Range("E3").Value = '  -0.15%  '

$ws.Range("E4").Value = '  -0.28%  '

$ws.Range("D5").Value = '218.81'
$ws.Range("E5").Value = '  -0.17%  '

$ws.Range("D6").Value = '0.5239'
$ws.Range("E6").Value = '  -0.21%  '

$ws.Range("E7").Value = '  -0.25%  '

$ws.Range("D8").Value = '0.2671'
$ws.Range("E8").Value = '  +1.42%  '

$ws.Range("D9").Value = '0.06360'
$ws.Range("E9").Value = '  +0.79%  '

$ws.Range("D10").Value = '20.63'
$ws.Range("E10").Value = '  -0.03%  '

$ws.Range("D11").Value = '0.07700'
$ws.Range("E11").Value = '  -1.42%  '

$ws.Range("D12").Value = '4.654'
$ws.Range("E12").Value = '  +3.61%  '

$ws.Range("D13").Value = '1.643.84'
$ws.Range("E13").Value = '  -1.64%  '

$ws.Range("D14").Value = '1.886.85'

$ws.Range("D15").Value = '0.5625'
$ws.Range("E15").Value = '  +1.40%  '

$ws.Range("D16").Value = '0.0₅8221'
$ws.Range("E16").Value = '  +2.41%  '

$ws.Range("D17").Value = '65.58'
$ws.Range("E17").Value = '  +0.60%  '

$ws.Range("D18").Value = '26.177.74'
$ws.Range("E18").Value = '  -0.16%  '

$ws.Range("E19").Value = '  -0.20%  '

$ws.Range("D20").Value = '4.670'
$ws.Range("E20").Value = '  +0.70%  '

$ws.Range("D21").Value = '10.49'
$ws.Range("E21").Value = '  +3.53%  '

$ws.Range("D22").Value = '192.26'

$ws.Range("D23").Value = '5.969'
$ws.Range("E23").Value = '  +0.05%  '

$ws.Range("E24").Value = '  -0.38%  '

$ws.Range("D25").Value = '145.74'
$ws.Range("E25").Value = '  +0.05%  '

$ws.Range("D26").Value = '0.1200'
$ws.Range("E26").Value = '  -0.45%  '

$ws.Range("D27").Value = '7.297'
$ws.Range("E27").Value = '  +2.08%  '

$ws.Range("D28").Value = '16.00'

$ws.Range("D29").Value = '1.526'
$ws.Range("E29").Value = '  +1.81%  '

$ws.Range("D30").Value = '0.05530'
$ws.Range("E30").Value = '  -4.07%  '

$ws.Range("E31").Value = '  -0.28%  '

$ws.Range("D32").Value = '3.479'
$ws.Range("E32").Value = '  -0.28%  '

$ws.Range("D33").Value = '3.385'
$ws.Range("E33").Value = '  +0.97%  '

$ws.Range("D34").Value = '1.569'
$ws.Range("E34").Value = '  -0.87%  '

$ws.Range("D35").Value = '0.9529'
$ws.Range("E35").Value = '  -0.12%  '

$ws.Range("D36").Value = '2.782'
$ws.Range("E36").Value = '  -0.84%  '

$ws.Range("D37").Value = '2.403'
$ws.Range("E37").Value = '  -0.73%  '

$ws.Range("D38").Value = '0.5718'
$ws.Range("E38").Value = '  +0.22%  '

$ws.Range("D39").Value = '0.01597'
$ws.Range("E39").Value = '  -0.03%  '

$ws.Range("D40").Value = '5.916'
$ws.Range("E40").Value = '  -0.75%  '

$ws.Range("E41").Value = '  -0.25%  '

$ws.Range("D42").Value = '1.033.94'
$ws.Range("E42").Value = '  -2.52%  '

$ws.Range("D43").Value = '0.8346'
$ws.Range("E43").Value = '  -1.90%  '

$ws.Range("D44").Value = '101.04'
$ws.Range("E44").Value = '  -2.22%  '

$ws.Range("D45").Value = '1.797.04'
$ws.Range("E45").Value = '  -0.13%  '

$ws.Range("D46").Value = '58.42'
$ws.Range("E46").Value = '  +0.41%  '

$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").Value = '0.9992'
$ws.Range("E47").Value = '  -1.12%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '8.056'
$ws.Range("E48").Value = '  +0.85%  '

$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").Value = '0.4347'
$ws.Range("E49").Value = '  -1.39%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.05242'
$ws.Range("E50").Value = '  +0.56%  '

$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₈101'
$ws.Range("E51").Value = '  +1.42%  '
